# Projet Courchevel - Soutenance finale : mise a jour des diagrammes de cas
# d'utilisation "Visiteur" (diapo 14) et "Accrediteur" (diapo 15), renommage
# d'une ellipse, et mise a jour des images illustrant les diapos 16-18
# (transparence du fond blanc + repositionnement/recadrage).

function EMU($v) {
    # Convertit des EMU (unite native OOXML) en points (unite COM PowerPoint).
    return $v / 12700.0
}

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# =======================================================================
# Diapo 14 : "Cas d'utilisation : Visiteur"
# =======================================================================
$s14 = $p.Slides.Item(14)

# Ellipse 11 (id 12) "Faire une demande individuelle presse" : deplacement
$shp = Get-ShapeById $s14 12
$shp.Left = EMU(3538358)
$shp.Top  = EMU(3861048)

# Connecteur droit 14 (id 15) : Picture 2 -> Ellipse 11 (id 12)
$shp = Get-ShapeById $s14 15
$shp.Left   = EMU(1965623)
$shp.Top    = EMU(3721559)
$shp.Width  = EMU(1572735)
$shp.Height = EMU(724611)

# Ellipse 17 (id 18) "Faire une demande groupee presse" : deplacement
$shp = Get-ShapeById $s14 18
$shp.Left = EMU(3491880)
$shp.Top  = EMU(5169390)

# Connecteur droit 27 (id 28) : Picture 2 -> Ellipse 17 (id 18)
$shp = Get-ShapeById $s14 28
$shp.Left   = EMU(1965623)
$shp.Top    = EMU(3721559)
$shp.Width  = EMU(1526257)
$shp.Height = EMU(1951887)

# =======================================================================
# Diapo 15 : "Cas d'utilisation : Accrediteur"
# =======================================================================
$s15 = $p.Slides.Item(15)

# Ellipse 9 (id 10) "Creer une accreditation individuelle" : deplacement
$shp = Get-ShapeById $s15 10
$shp.Left = EMU(3779912)
$shp.Top  = EMU(1960180)

# Ellipse 10 (id 11) "Creer une accreditation groupee" : deplacement
$shp = Get-ShapeById $s15 11
$shp.Left = EMU(4355976)
$shp.Top  = EMU(3239438)

# Ellipse 11 (id 12) "Valider les accreditations" : deplacement
$shp = Get-ShapeById $s15 12
$shp.Left = EMU(3779912)
$shp.Top  = EMU(4437112)

# Connecteur droit 12 (id 13) : Picture 2 -> Ellipse 9 (id 10) ; suit le
# deplacement de l'ellipse 9 et devient retourne verticalement.
$shp = Get-ShapeById $s15 13
$shp.Left          = EMU(1349608)
$shp.Top           = EMU(2530235)
$shp.Width         = EMU(2430304)
$shp.Height        = EMU(1086648)
$shp.VerticalFlip  = -1

# Connecteur droit 13 (id 14) : Picture 2 -> Ellipse 10 (id 11) ; suit le
# deplacement de l'ellipse 10.
$shp = Get-ShapeById $s15 14
$shp.Left   = EMU(1349608)
$shp.Top    = EMU(3616883)
$shp.Width  = EMU(3006368)
$shp.Height = EMU(162615)

# Connecteur droit 14 (id 15) : Picture 2 -> Ellipse 11 (id 12) ; suit le
# deplacement de l'ellipse 11 et n'est plus retourne verticalement.
$shp = Get-ShapeById $s15 15
$shp.Left          = EMU(1349608)
$shp.Top           = EMU(3616883)
$shp.Width         = EMU(2430304)
$shp.Height        = EMU(1405351)
$shp.VerticalFlip  = 0

# Ellipse 52 (id 53) : le texte "Faire une demande groupee presse" devient
# "Imprimer une accreditation".
$shp = Get-ShapeById $s15 53
$tr = $shp.TextFrame.TextRange
$sub = $tr.Characters(19, 14)
$sub.Text = "Imprimer une accréditation"
$del = $tr.Characters(1, 18)
$del.Text = ""

# =======================================================================
# Diapos 16-18 : images "Picture 2" -> fond blanc transparent + recadrage
# =======================================================================

# Diapo 16
$s16 = $p.Slides.Item(16)
$shp = Get-ShapeById $s16 2050
$shp.PictureFormat.TransparencyColor = 16777215
$shp.PictureFormat.TransparentBackground = -1
$shp.Left   = EMU(432000)
$shp.Top    = EMU(927056)
$shp.Width  = EMU(8542240)
$shp.Height = EMU(5580000)

# Diapo 17
$s17 = $p.Slides.Item(17)
$shp = Get-ShapeById $s17 3074
$shp.PictureFormat.TransparencyColor = 16777215
$shp.PictureFormat.TransparentBackground = -1
$shp.Left   = EMU(432000)
$shp.Top    = EMU(928800)
$shp.Width  = EMU(8542237)
$shp.Height = EMU(5580000)

# Diapo 18
$s18 = $p.Slides.Item(18)
$shp = Get-ShapeById $s18 4098
$shp.PictureFormat.TransparencyColor = 16777215
$shp.PictureFormat.TransparentBackground = -1
$shp.Left   = EMU(432480)
$shp.Top    = EMU(927054)
$shp.Width  = EMU(8542237)
$shp.Height = EMU(5580000)

Write-Output "Done."
